$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Z2:Z10").Value = "2025-11-02T02:03:33.550474"
$ws.Range("Z11:Z20").Value = "2025-11-02T02:03:33.551474"
$ws.Range("Z21").Value = "2025-11-02T02:03:33.552473"
$ws.Range("Z22:Z25").Value = "2025-11-02T02:03:33.552582"
$ws.Range("Z26:Z29").Value = "2025-11-02T02:03:33.553114"
